# Update cryptos list prices/volumes (and the Polygon/WrappedEther row swap)
# to match the latest scrape, per the commit "Updated cryptos list ...".
#
# Column D ("Price") cells that look like plain numbers (e.g. "4.57") are
# forced to stay text the same way the source data does, by temporarily
# marking the cell as Text before assignment and then clearing the
# resulting formatting so the cell's style index is left untouched
# (matching the original file, where these cells carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.678.19'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.809.20'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.16'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -10.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.324'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.99%  '
$ws.Range("E10").Value = '  -3.08%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.072.78'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.810.52'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.663'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.57'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.674.35'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.30'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.84'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.85'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.66'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("E24").Value = '  +2.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.50'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.07'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +16.33%  '
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.696'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '91.51'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.00%  '
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.326.93'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0191'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("E41").Value = '  +0.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.958'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.27'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -8.04%  '
$ws.Range("E44").Value = '  -9.17%  '
$ws.Range("E45").Value = '  -4.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.23'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0514'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.999.47'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  +5.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.71'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.06%  '
